$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column K (11) to hold the "ImagingMethod" field,
# shifting the existing AdMethod..ScopeValueName columns one to the right.
$ws.Columns.Item(11).Insert()

# New header + value for the inserted "ImagingMethod" column.
$ws.Cells.Item(1, 11).Value = "ImagingMethod"
$ws.Cells.Item(2, 11).Value = "Absorption"
$ws.Columns.Item(11).ColumnWidth = 16.285714285714285

$ws.Range("K3").Select()
